$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G4").Value = "2017-02-21 02:56:59"

$wsZhCn.Range("H4").Value = "2017-02-21 02:56:43"
$wsZhCn.Range("L4").Value = "2017-02-21 02:57:40"

$wsDeDe.Range("H4").Value = "2017-02-21 02:56:59"
$wsDeDe.Range("L4").Value = "2017-02-21 02:58:02"
